# Testy manualne - add row 9 ("Koszyk" test case) to the first worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row content ------------------------------------------------------
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Koszyk"
$ws.Range("C9").Value = "Szczegół"
$ws.Range("I9").Value = "po kliknięciu na koszyk bez check-a "
$ws.Range("G9").Value = "L"

# Hyperlink for D9 (adds the relationship + shared string "Wordy\Koszyk.docx").
$null = $ws.Hyperlinks.Add($ws.Range("D9"), "Wordy\Koszyk.docx")
# Re-apply the same "hyperlink" cell style used by the other link cells
# (D2/D4/D5/D8) instead of the fresh style Hyperlinks.Add auto-generates.
$ws.Range("D9").Style = $ws.Range("D8").Style

# Row 9 mirrors the other multi-line rows (1,2,4,5,6) at 28.5pt.
$ws.Rows.Item(9).RowHeight = 28.5

# Move the active selection to F9, matching the saved view state.
$null = $ws.Range("F9").Select()
